# Apply the targeted updates to column F (dSF) as described by the diff.
# Only specific cells in column F change value; everything else is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = 4
    11 = 4
    13 = -1
    17 = 0
    21 = 1
    22 = -1
    23 = -1
    25 = 1
    29 = -1
    32 = 0
    34 = -1
    35 = 4
    37 = -1
    41 = -1
    51 = -3
    54 = 0
    55 = 11
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
